$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 624
$ws.Range("F3").Value = 5898
$ws.Range("F8").Value = 398
$ws.Range("F11").Value = 3125
$ws.Range("F12").Value = 1961
$ws.Range("F15").Value = 202
$ws.Range("F16").Value = 83
$ws.Range("F17").Value = 177
$ws.Range("F19").Value = 996
$ws.Range("F20").Value = 366
$ws.Range("F21").Value = 59
$ws.Range("F22").Value = 65
$ws.Range("F23").Value = 3682
$ws.Range("F24").Value = 1170
$ws.Range("F25").Value = 2924
$ws.Range("F27").Value = 2352
$ws.Range("F28").Value = 4244
$ws.Range("F29").Value = 114
$ws.Range("F30").Value = 933
$ws.Range("F31").Value = 481
$ws.Range("F33").Value = 98
$ws.Range("F35").Value = 31
$ws.Range("F36").Value = 29
$ws.Range("F38").Value = 1030
$ws.Range("F39").Value = 1292
$ws.Range("F40").Value = 76
$ws.Range("F41").Value = 1107
$ws.Range("F42").Value = 717
$ws.Range("F43").Value = 610
$ws.Range("F44").Value = 434
$ws.Range("F45").Value = 17
$ws.Range("F46").Value = 109
$ws.Range("F49").Value = 3612

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 10
$ws.Range("F23").Value = 40
$ws.Range("F25").Value = 19

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 624
$ws.Range("F3").Value = 5898
$ws.Range("F6").Value = 10
$ws.Range("F7").Value = 398
$ws.Range("F9").Value = 3125
$ws.Range("F11").Value = 1961
$ws.Range("F14").Value = 202
$ws.Range("F17").Value = 83
$ws.Range("F18").Value = 177
$ws.Range("F19").Value = 996
$ws.Range("F20").Value = 366
$ws.Range("F21").Value = 3682
$ws.Range("F23").Value = 1170
$ws.Range("F25").Value = 2924
$ws.Range("F26").Value = 2355
$ws.Range("F27").Value = 4244
$ws.Range("F28").Value = 114
$ws.Range("F29").Value = 933
$ws.Range("F32").Value = 1030
$ws.Range("F34").Value = 1292
$ws.Range("F35").Value = 76
$ws.Range("F36").Value = 1107
$ws.Range("F38").Value = 717
$ws.Range("F40").Value = 434
$ws.Range("F41").Value = 40
$ws.Range("F43").Value = 17
$ws.Range("F44").Value = 19
$ws.Range("F45").Value = 109
$ws.Range("F48").Value = 3612
